$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.65126913244769
$ws.Range("D2").Value = 5.438758906836134
$ws.Range("E2").Value = 16.45815516714121
$ws.Range("F2").Value = 28.57045719014965
$ws.Range("G2").Value = 36.44697243724561
$ws.Range("H2").Value = 15.96245578658778
$ws.Range("K2").Value = 9.536744868163485
$ws.Range("L2").Value = 8.828641205660858
$ws.Range("N2").Value = 20.57779678225955
$ws.Range("B3").Value = 16.54798198233219
$ws.Range("D3").Value = 5.429752380233798
$ws.Range("E3").Value = 16.49709415470262
$ws.Range("F3").Value = 28.485987986974
$ws.Range("G3").Value = 36.26388192742439
$ws.Range("H3").Value = 15.98575378650728
$ws.Range("K3").Value = 9.227148350241642
$ws.Range("L3").Value = 8.789046906251006
$ws.Range("N3").Value = 20.64304890491348
$ws.Range("B4").Value = 16.48798100918538
$ws.Range("D4").Value = 5.424112203045601
$ws.Range("E4").Value = 16.52255671698286
$ws.Range("F4").Value = 28.44193479066537
$ws.Range("G4").Value = 36.16268690329398
$ws.Range("H4").Value = 16.00312481525255
$ws.Range("K4").Value = 9.033873482264488
$ws.Range("L4").Value = 8.766455163055239
$ws.Range("N4").Value = 20.68501532554269
$ws.Range("B5").Value = 16.46441142218822
$ws.Range("D5").Value = 5.4217864365063
$ws.Range("E5").Value = 16.53332439394014
$ws.Range("F5").Value = 28.42595863949198
$ws.Range("G5").Value = 36.12430348909569
$ws.Range("H5").Value = 16.0109735000431
$ws.Range("K5").Value = 8.954445479570687
$ws.Range("L5").Value = 8.757688305313113
$ws.Range("N5").Value = 20.70259632894252
$ws.Range("B6").Value = 16.46055154986297
$ws.Range("D6").Value = 5.421398603088636
$ws.Range("E6").Value = 16.53513602464521
$ws.Range("F6").Value = 28.42342544427806
$ws.Range("G6").Value = 36.11810319530476
$ws.Range("H6").Value = 16.01232322916491
$ws.Range("K6").Value = 8.941220322383083
$ws.Range("L6").Value = 8.756259314693695
$ws.Range("N6").Value = 20.70554462889917
$ws.Range("B7").Value = 16.4876595453145
$ws.Range("D7").Value = 5.424080947231503
$ws.Range("E7").Value = 16.52270034756041
$ws.Range("F7").Value = 28.44171131684118
$ws.Range("G7").Value = 36.16215765502395
$ws.Range("H7").Value = 16.00322754990304
$ws.Range("K7").Value = 9.032804805679376
$ws.Range("L7").Value = 8.766335141794443
$ws.Range("N7").Value = 20.68525048662637
$ws.Range("B8").Value = 16.61495991292978
$ws.Range("D8").Value = 5.435676519664994
$ws.Range("E8").Value = 16.47125927249154
$ws.Range("F8").Value = 28.53971872096567
$ws.Range("G8").Value = 36.38153421588855
$ws.Range("H8").Value = 15.96985197829855
$ws.Range("K8").Value = 9.430732993890773
$ws.Range("L8").Value = 8.814636689924797
$ws.Range("N8").Value = 20.5999018352999
$ws.Range("B9").Value = 16.8906889866595
$ws.Range("D9").Value = 5.457527115318526
$ws.Range("E9").Value = 16.38267916451618
$ws.Range("F9").Value = 28.79326663798089
$ws.Range("G9").Value = 36.89914091246938
$ws.Range("H9").Value = 15.92877332308687
$ws.Range("K9").Value = 10.18049775875078
$ws.Range("L9").Value = 8.922659759059465
$ws.Range("N9").Value = 20.44756227663083
$ws.Range("B10").Value = 17.10771838304982
$ws.Range("D10").Value = 5.473020949936123
$ws.Range("E10").Value = 16.32504873862548
$ws.Range("F10").Value = 29.01596255757254
$ws.Range("G10").Value = 37.33020391111756
$ws.Range("H10").Value = 15.91350063921828
$ws.Range("K10").Value = 10.83672947963153
$ws.Range("L10").Value = 9.009647549778879
$ws.Range("N10").Value = 20.34471834831893
$ws.Range("B11").Value = 17.20925709143772
$ws.Range("D11").Value = 5.47994413373819
$ws.Range("E11").Value = 16.30043896994118
$ws.Range("F11").Value = 29.12492107051287
$ws.Range("G11").Value = 37.5367095963224
$ws.Range("H11").Value = 15.90979664412278
$ws.Range("K11").Value = 11.13871659253766
$ws.Range("L11").Value = 9.050762224664156
$ws.Range("N11").Value = 20.2998861548747
$ws.Range("B12").Value = 17.24808226087215
$ws.Range("D12").Value = 5.4825474396428
$ws.Range("E12").Value = 16.29135018823516
$ws.Range("F12").Value = 29.1672561211228
$ws.Range("G12").Value = 37.6163474467814
$ws.Range("H12").Value = 15.90886062376867
$ws.Range("K12").Value = 11.25070878917686
$ws.Range("L12").Value = 9.066543007690965
$ws.Range("N12").Value = 20.28318880254322
$ws.Range("B13").Value = 17.23970437947602
$ws.Range("D13").Value = 5.481987593334237
$ws.Range("E13").Value = 16.29329738257934
$ws.Range("F13").Value = 29.15809113735688
$ws.Range("G13").Value = 37.59913299637505
$ws.Range("H13").Value = 15.90904145988342
$ws.Range("K13").Value = 11.22669471555187
$ws.Range("L13").Value = 9.063135094497106
$ws.Range("N13").Value = 20.28677245657889
$ws.Range("B14").Value = 17.21244392404488
$ws.Range("D14").Value = 5.48015867924576
$ws.Range("E14").Value = 16.29968661630998
$ws.Range("F14").Value = 29.12838261366253
$ws.Range("G14").Value = 37.54323296572114
$ws.Range("H14").Value = 15.90971028557557
$ws.Range("K14").Value = 11.14797774804778
$ws.Range("L14").Value = 9.052056340306178
$ws.Range("N14").Value = 20.29850685561216
$ws.Range("B15").Value = 17.19579401300373
$ws.Range("D15").Value = 5.47903601533984
$ws.Range("E15").Value = 16.30363019344263
$ws.Range("F15").Value = 29.11032446295639
$ws.Range("G15").Value = 37.50917813682463
$ws.Range("H15").Value = 15.91018072728171
$ws.Range("K15").Value = 11.09945288897283
$ws.Range("L15").Value = 9.045297518355341
$ws.Range("N15").Value = 20.30573089786891
$ws.Range("B16").Value = 17.10113664182663
$ws.Range("D16").Value = 5.472565955553837
$ws.Range("E16").Value = 16.32668932029886
$ws.Range("F16").Value = 29.00899371316782
$ws.Range("G16").Value = 37.31691286735857
$ws.Range("H16").Value = 15.91380799394959
$ws.Range("K16").Value = 10.8166643011958
$ws.Range("L16").Value = 9.006990787048046
$ws.Range("N16").Value = 20.34768742531016
$ws.Range("B17").Value = 17.04376654601309
$ws.Range("D17").Value = 5.468564491933908
$ws.Range("E17").Value = 16.34124639032095
$ws.Range("F17").Value = 28.94877309165365
$ws.Range("G17").Value = 37.20159061684542
$ws.Range("H17").Value = 15.91686416907262
$ws.Range("K17").Value = 10.63899606653514
$ws.Range("L17").Value = 8.983878888741566
$ws.Range("N17").Value = 20.37392558860695
$ws.Range("B18").Value = 17.0110351775851
$ws.Range("D18").Value = 5.466251208772249
$ws.Range("E18").Value = 16.34977048967507
$ws.Range("F18").Value = 28.91485779160904
$ws.Range("G18").Value = 37.13624416249468
$ws.Range("H18").Value = 15.9189273092436
$ws.Range("K18").Value = 10.53527532960426
$ws.Range("L18").Value = 8.970731464372586
$ws.Range("N18").Value = 20.38920088312304
$ws.Range("B19").Value = 16.99999954170276
$ws.Range("D19").Value = 5.465465964667809
$ws.Range("E19").Value = 16.3526826010383
$ws.Range("F19").Value = 28.90349940809259
$ws.Range("G19").Value = 37.11428968984581
$ws.Range("H19").Value = 15.9196782812764
$ws.Range("K19").Value = 10.49989519118464
$ws.Range("L19").Value = 8.966305346683361
$ws.Range("N19").Value = 20.39440443563968
$ws.Range("B20").Value = 17.04984633197579
$ws.Range("D20").Value = 5.468991673047956
$ws.Range("E20").Value = 16.33968111532128
$ws.Range("F20").Value = 28.95510913386614
$ws.Range("G20").Value = 37.2137654572448
$ws.Range("H20").Value = 15.91650723367457
$ws.Range("K20").Value = 10.65806776102458
$ws.Range("L20").Value = 8.986324160958135
$ws.Range("N20").Value = 20.37111347682887
$ws.Range("B21").Value = 17.22044105078551
$ws.Range("D21").Value = 5.48069637672114
$ws.Range("E21").Value = 16.29780369532764
$ws.Range("F21").Value = 29.13707977127246
$ws.Range("G21").Value = 37.55961359969949
$ws.Range("H21").Value = 15.90950117154664
$ws.Range("K21").Value = 11.17116314375411
$ws.Range("L21").Value = 9.055304782297931
$ws.Range("N21").Value = 20.29505259766348
$ws.Range("B22").Value = 17.33410429413605
$ws.Range("D22").Value = 5.488238926841927
$ws.Range("E22").Value = 16.27177706561831
$ws.Range("F22").Value = 29.26225971679058
$ws.Range("G22").Value = 37.79400137923619
$ws.Range("H22").Value = 15.90764193980083
$ws.Range("K22").Value = 11.49271542642806
$ws.Range("L22").Value = 9.101615602023799
$ws.Range("N22").Value = 20.2469718054854
$ws.Range("B23").Value = 17.27325146064988
$ws.Range("D23").Value = 5.484223252969789
$ws.Range("E23").Value = 16.28554531576589
$ws.Range("F23").Value = 29.19488579398921
$ws.Range("G23").Value = 37.66815984242497
$ws.Range("H23").Value = 15.90838540141488
$ws.Range("K23").Value = 11.32236474447205
$ws.Range("L23").Value = 9.076789811223993
$ws.Range("N23").Value = 20.27248469042944
$ws.Range("B24").Value = 17.04709687520403
$ws.Range("D24").Value = 5.468798584252629
$ws.Range("E24").Value = 16.34038829317384
$ws.Range("F24").Value = 28.95224240699633
$ws.Range("G24").Value = 37.2082582362343
$ws.Range("H24").Value = 15.91666765063982
$ws.Range("K24").Value = 10.64945035788417
$ws.Range("L24").Value = 8.985218216787324
$ws.Range("N24").Value = 20.37238423865547
$ws.Range("B25").Value = 16.81344942890593
$ws.Range("D25").Value = 5.451712939403131
$ws.Range("E25").Value = 16.40533092143595
$ws.Range("F25").Value = 28.71819960144747
$ws.Range("G25").Value = 36.74999881197252
$ws.Range("H25").Value = 15.9372710888194
$ws.Range("K25").Value = 9.981640757603133
$ws.Range("L25").Value = 8.892060711459147
$ws.Range("N25").Value = 20.48717375456373
